# Commit: "se agrega 6097 y 6089" -- add two new part rows (2098706087 / 2098706089)
# to the "Hoja1" parts table, which is kept sorted by column A (part number).
# Excel's Table auto-sort places them right before 2099700009 (row 99),
# pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Make room for the two new records by inserting two blank rows at the
# position where they belong in sorted order (just above the old row 99).
$ws.Range("A99:A100").EntireRow.Insert()

# Fill in the two new rows. Column A keeps the "0" number format that the
# rest of the table's A column intermittently carries (style index 5).
$ws.Range("A99").Value = 2098706087
$ws.Range("A99").NumberFormat = "0"
$ws.Range("B99").Value = "V1.1"
$ws.Range("C99").Value = "AMZW17-000-C"
$ws.Range("D99").Value = "59Z118-C00-E"

$ws.Range("A100").Value = 2098706089
$ws.Range("A100").NumberFormat = "0"
$ws.Range("B100").Value = "V1.1"
$ws.Range("C100").Value = "AMZW17-000-C"
$ws.Range("D100").Value = "59Z118-C00-E"

# The Table ("Tabla1") covers this data; grow it to include the two new
# rows and re-apply its sort (by Col 1 / column A) so the table + worksheet
# sort metadata stay consistent with the new A1:H230 extent.
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.Resize($ws.Range("A1:H230"))

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A230"))
$sortObj.SetRange($ws.Range("A2:G230"))
$sortObj.Header = 0
$sortObj.Apply()

# Mirror the author's final on-screen state: scrolled/selected near the new
# rows on Hoja1, then ended up on Hoja3 ("Hoja oficial de busca" stays put).
$ws.Range("C100:D100").Select()

$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Select()
$ws3.Range("AS6").Select()
